# Download-Format-Item.xlsx — Inward QC / purchase-order item-master sample rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New second sample item (Sell / Semi-Finished Good / TM_0002) entered first
$ws.Range("B3").Value = "Sell"
$ws.Range("C3").Value = "Semi-Finished Good"
$ws.Range("D3").Value = "TM_0002"

# Rename the sample item names in column A
$ws.Range("A2").Value = "ABC"
$ws.Range("A3").Value = "XYZ"

# Row 2 remaining columns (kept same as before: Buy / Raw Material / TM_0001 / 10000)
$ws.Range("B2").Value = "Buy"
$ws.Range("C2").Value = "Raw Material"
$ws.Range("D2").Value = "TM_0001"
$ws.Range("F2").Value = 10000

# Row 3 stock value
$ws.Range("F3").Value = 3000

# Restore the selection/active cell as left by the author after editing
[void]$ws.Range("C11").Select()
